$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the KODE_REKSADANA value in N2 from "RD00014" to "RD00015"
$ws.Range("N2").Value = "RD00015"

# Update the PREPARATION text in F2 to reference the new code
$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 18/19/20/21 - Pimpinan Kelompok Investasi/Pengelolan Investasi/Analis;`nKode Reksadana : RD00015;`nNama Reksadana : Reksadana BNI-QAO Syariah"

# Update the selected cell in the sheet view from O2 to G2
$ws.Range("G2").Select()
